$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 477, shifting existing rows 477-532 down to 478-533
$ws.Rows.Item(477).Insert()

# Populate the newly inserted row 477 with the new data point
$ws.Cells.Item(477, 1).Value = 10
$ws.Cells.Item(477, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(477, 3).Value = "La Araucanía"
$ws.Cells.Item(477, 4).Value = 44858
$ws.Cells.Item(477, 5).Value = 9
$ws.Cells.Item(477, 6).Value = 100112023
$ws.Cells.Item(477, 7).Value = "Brócoli"
$ws.Cells.Item(477, 8).Value = "Sin especificar"
$ws.Cells.Item(477, 9).Value = "Primera"
$ws.Cells.Item(477, 10).Value = 2500
$ws.Cells.Item(477, 11).Value = 1000
$ws.Cells.Item(477, 12).Value = 1000
$ws.Cells.Item(477, 13).Value = 1000
$ws.Cells.Item(477, 14).Value = "$/unidad"
$ws.Cells.Item(477, 15).Value = "Región del Maule"
$ws.Cells.Item(477, 16).Value = 1000
$ws.Cells.Item(477, 17).Value = 1
$ws.Cells.Item(477, 18).Value = "Hortaliza"
